$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7, shifting the current row 7 (and its
# formatting) down to row 8. This preserves the D8 style (date format)
# that was previously on row 7.
$ws.Rows.Item(7).Insert()

# New row 7 values (new weekly data point)
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(7, 4).Value = 44798
$ws.Cells.Item(7, 5).Value = 15
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100108
$ws.Cells.Item(7, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(7, 9).Value = 100108007
$ws.Cells.Item(7, 10).Value = "Coco"
$ws.Cells.Item(7, 11).Value = "Sin especificar"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 80
$ws.Cells.Item(7, 14).Value = 21000
$ws.Cells.Item(7, 15).Value = 22000
$ws.Cells.Item(7, 16).Value = 21500
$ws.Cells.Item(7, 17).Value = "`$/malla 20 unidades"
$ws.Cells.Item(7, 18).Value = "Perú"
$ws.Cells.Item(7, 19).Value = 1075
$ws.Cells.Item(7, 20).Value = 20

# Ensure row 7's date cell (D7) uses the same date style as before.
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat

# Row 8 keeps the original row-7 data (shifted down by the insert), but
# set the values explicitly to make sure they match exactly.
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value = 44533
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100108
$ws.Cells.Item(8, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(8, 9).Value = 100108007
$ws.Cells.Item(8, 10).Value = "Coco"
$ws.Cells.Item(8, 11).Value = "Sin especificar"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 100
$ws.Cells.Item(8, 14).Value = 16000
$ws.Cells.Item(8, 15).Value = 17000
$ws.Cells.Item(8, 16).Value = 16500
$ws.Cells.Item(8, 17).Value = "`$/malla 20 unidades"
$ws.Cells.Item(8, 18).Value = "Perú"
$ws.Cells.Item(8, 19).Value = 825
$ws.Cells.Item(8, 20).Value = 20
